$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '63.713.24'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -6.45%  '

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.266.11'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -9.00%  '

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.00'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  +0.12%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '177.11'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -12.63%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '510.33'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -10.15%  '

# Row 7
$ws.Range("B7").Value = 'XRP'
$ws.Range("C7").Value = 'https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp'
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.584'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -4.61%  '

# Row 8
$ws.Range("B8").Value = 'LidoStakedEther'
$ws.Range("C8").Value = 'https://coinranking.com/coin/VINVMYf0u+lidostakedether-steth'
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '3.266.29'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -8.97%  '

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '1.00'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +0.07%  '

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.609'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -10.28%  '

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '57.21'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -5.03%  '

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.129'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -12.48%  '

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.0000250'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -10.81%  '

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '8.93'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -12.58%  '

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '3.774.28'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -9.00%  '

# Row 16
$ws.Range("E16").Value = '  -4.92%  '

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '3.258.82'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -8.88%  '

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '63.434.11'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -6.41%  '

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '16.89'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -11.09%  '

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '10.67'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -12.81%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.932'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -12.03%  '

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '365.74'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -9.19%  '

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '79.00'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -6.68%  '

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '3.61'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -13.47%  '

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '10.67'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -14.74%  '

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '5.96'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -2.91%  '

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '3.74'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -3.86%  '

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.59'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -10.06%  '

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '11.09'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -10.54%  '

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '8.21'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -10.69%  '

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '28.06'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -10.91%  '

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '634.89'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -5.52%  '

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '6.57'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -15.06%  '

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '10.93'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -9.57%  '

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '58.50'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -7.61%  '

# Row 36
$ws.Range("B36").Value = 'Hedera'
$ws.Range("C36").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.101'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -9.86%  '

# Row 37
$ws.Range("B37").Value = 'Dai'
$ws.Range("C37").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.00'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -0.05%  '

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '35.20'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -14.26%  '

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.368'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -9.59%  '

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.996'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +0.02%  '

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.121'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -9.20%  '

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '2.806.00'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -11.89%  '

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '2.65'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -18.38%  '

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.0₃0625'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -17.47%  '

# Row 45
$ws.Range("B45").Value = 'WEMIXToken'
$ws.Range("C45").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '2.56'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -8.96%  '

# Row 46
$ws.Range("B46").Value = 'VeChain'
$ws.Range("C46").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.0377'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -8.12%  '

# Row 47
$ws.Range("B47").Value = 'Fetch.AI'
$ws.Range("C47").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.25'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -15.50%  '

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.122'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -6.38%  '

# Row 49
$ws.Range("B49").Value = 'Monero'
$ws.Range("C49").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '131.70'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -5.39%  '

# Row 50
$ws.Range("B50").Value = 'Stacks'
$ws.Range("C50").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '2.63'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -2.51%  '

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '2.72'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -11.36%  '

